$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert the four new "XGBoost" rows. Working from the bottom of the
#    sheet upward so that each insertion point is still expressed in the
#    sheet's ORIGINAL row numbering (rows below each insertion point shift
#    down by one, rows above are untouched).
# ---------------------------------------------------------------------------
$ws.Rows.Item(53).Insert()   # "Si" block: before Logistic Regression
$ws.Rows.Item(31).Insert()   # "C, Si, N, Al" block: before Random Forest Classifier
$ws.Rows.Item(24).Insert()   # "CARBON and SILICON" block: before Random Forest Classifier
$ws.Rows.Item(11).Insert()   # "CARBON ONLY" block: before kNN

# ---------------------------------------------------------------------------
# 2. Fix up the title text (lower-case "comparison" -> "Comparison").
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Machine Learning Model Comparison for Meteorite Classification"

# ---------------------------------------------------------------------------
# 3. Add the "# records" annotation cells to the three title rows that were
#    missing them (CARBON ONLY / CARBON and SILICON / C, Si, N, Al).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = 15622
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("D2").Value = "# records"

$ws.Range("C19").Value = 14423
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("D19").Value = "# records"

$ws.Range("C28").Value = 328
$ws.Range("D28").Value = "# records"

# ---------------------------------------------------------------------------
# 4. Populate the newly-inserted "XGBoost" rows with their model name/rating.
# ---------------------------------------------------------------------------
$ws.Range("A11").Value = "XGBoost"
$ws.Range("B11").Value = 93.7

$ws.Range("A25").Value = "XGBoost"
$ws.Range("B25").Value = 96.8

$ws.Range("A33").Value = "XGBoost"
$ws.Range("B33").Value = 95.1

$ws.Range("A56").Value = "XGBoost"
$ws.Range("B56").Value = 87.8

# ---------------------------------------------------------------------------
# 5. The existing sort annotation on the "C, Si, N" block ("A42:B47") needs
#    to be re-anchored to its new location ("A45:B50") now that four rows
#    were inserted above it.
# ---------------------------------------------------------------------------
$ws.Sort.SetRange($ws.Range("A45:B50"))
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B45:B50"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 6. Misc view/selection tidy-up to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("B1").Select()

Write-Output "done"
